$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has data rows 2..855 (used range A1:R855).
# This edit inserts two brand-new data rows right before the existing row 819
# (shifting the old rows 819..855 down to 821..857, preserving all of their
# values/formatting), and populates the two new rows (819, 820) with fresh
# data. Net effect: dimension grows from A1:R855 to A1:R857.

$ws.Rows.Item(819).Resize(2).Insert()

# ---- New row 819 ----
$ws.Range("A819").Value = 4
$ws.Range("B819").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C819").Value = "Los Lagos"
$ws.Range("D819").Value = 45041
$ws.Range("E819").Value = 10
$ws.Range("F819").Value = 100112033
$ws.Range("G819").Value = "Lechuga"
$ws.Range("H819").Value = "Escarola"
$ws.Range("I819").Value = "Primera"
$ws.Range("J819").Value = 300
$ws.Range("K819").Value = 12000
$ws.Range("L819").Value = 12000
$ws.Range("M819").Value = 12000
$ws.Range("N819").Value = "`$/caja 15 unidades"
$ws.Range("O819").Value = "Región de Coquimbo"
$ws.Range("P819").Value = 800
$ws.Range("Q819").Value = 15
$ws.Range("R819").Value = "Hortaliza"

# ---- New row 820 ----
$ws.Range("A820").Value = 4
$ws.Range("B820").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C820").Value = "Los Lagos"
$ws.Range("D820").Value = 45041
$ws.Range("E820").Value = 10
$ws.Range("F820").Value = 100112033
$ws.Range("G820").Value = "Lechuga"
$ws.Range("H820").Value = "Escarola"
$ws.Range("I820").Value = "Segunda"
$ws.Range("J820").Value = 300
$ws.Range("K820").Value = 10500
$ws.Range("L820").Value = 10500
$ws.Range("M820").Value = 10500
$ws.Range("N820").Value = "`$/caja 18 unidades"
$ws.Range("O820").Value = "Región de Coquimbo"
$ws.Range("P820").Value = 583
$ws.Range("Q820").Value = 18
$ws.Range("R820").Value = "Hortaliza"
